$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the three runs "a publicidad buzoneando" / " el siguiente flayer" /
#    " por las zonas Triana, " into a single run, leaving the surrounding
#    runs ("...que realizar") and ("Santa Justa"...) untouched.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*buzoneando*") {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EF1DDB" w:rsidRDefault="00EF1DDB" w:rsidP="00D12ACA"><w:r><w:t xml:space="preserve">También contamos con la participación de SVQ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Fixie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, un grupo local de Sevilla amante de las bicicletas de piñón fijo que realizar</w:t></w:r><w:r><w:t xml:space="preserve">a publicidad buzoneando el siguiente flayer por las zonas Triana, </w:t></w:r><w:r><w:t>Santa Justa</w:t></w:r><w:r w:rsidR="007E39AA"><w:t xml:space="preserve"> y La Alameda de Hércules.</w:t></w:r></w:p>'
        $p.Range.InsertXML($xml) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Grow the API routes table: the last row only held the "_GoBack" bookmark
#    (empty route/verb/action). Turn it into three fully populated rows:
#       /profile/{id}  PUT     Edita el perfil del usuario
#       /pub           POST    Crea una nueva publicación
#       /pub/{id}      DELETE  Elimina una publicación del usuario por el id de la publicacion
#    and keep the _GoBack bookmark, relocated to the end of the last cell.
# ---------------------------------------------------------------------------
$tbl = $null
foreach ($t in $d.Tables) {
    $headerText = $t.Rows.Item(1).Cells.Item(1).Range.Text
    if ($headerText -like "*Ruta*") {
        $tbl = $t
        break
    }
}

$origLastIndex = $tbl.Rows.Count

$rowA = $tbl.Rows.Add()
$rowB = $tbl.Rows.Add()
$rowC = $tbl.Rows.Add()

# Row A: /profile/{id} | PUT | Edita el perfil del usuario
$xmlA1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>profile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/{id}</w:t></w:r></w:p>'
$rowA.Cells.Item(1).Range.Paragraphs.Item(1).Range.InsertXML($xmlA1) | Out-Null
$rowA.Cells.Item(2).Range.Text = "PUT"
$rowA.Cells.Item(3).Range.Text = "Edita el perfil del usuario"

# Row B: /pub | POST | Crea una nueva publicación
$rowB.Cells.Item(1).Range.Text = "/pub"
$rowB.Cells.Item(2).Range.Text = "POST"
$rowB.Cells.Item(3).Range.Text = "Crea una nueva publicación"

# Row C: /pub/{id} | DELETE | Elimina una publicación... + _GoBack bookmark
$rowC.Cells.Item(1).Range.Text = "/pub/{id}"
$rowC.Cells.Item(2).Range.Text = "DELETE"
$xmlC3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Elimina una publicación del usuario por el id de la publicacion</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$rowC.Cells.Item(3).Range.Paragraphs.Item(1).Range.InsertXML($xmlC3) | Out-Null

# Remove the old, now-superseded, bookmark-only row.
$tbl.Rows.Item($origLastIndex).Delete() | Out-Null

Write-Output "done"
